$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new label in J6 describing the "Most Healthy people don't get treatment" scenario
$ws.Range("J6").Value = "Most Healthy people don't get treatment"

# Update the "No treatment" CPT table (rows 8-9, columns C-D)
$ws.Range("C8").Value = 0.584
$ws.Range("D8").Value = 0.416
$ws.Range("C9").Value = 0.867
$ws.Range("D9").Value = 0.133

# Update the selected cell in the sheet view
$ws.Range("D12").Select()
